# Daily attendance processing - 2026-01-21 08:46:25
#
# This script applies the daily attendance-processing pass described in the
# commit:
#   1. Six sessions dated 21/01/2026 (one per B1 sub-group: B1D1, B1D2, B1E1,
#      B1E2, B1F1, B1F2) have now passed their scheduled time without being
#      recorded, so each one flips from "Pending" (yellow) to
#      "Not Recorded" (pink) - matching the formatting/status already used
#      for older un-recorded sessions.
#   2. The roll-up "Class Statistics" counters move accordingly: Missing
#      Sessions (L7) goes up by 6, Pending Sessions (L8) goes down by 6.
#   3. The per-group "Group Statistics" table (rows 21-26, groups
#      B1D1/B1D2/B1E1/B1E2/B1F1/B1F2) mirrors the same shift: Missing (P)
#      +1, Pending (Q) -1 for each of those six groups.
#   4. The "Recorded By" column for every previously-recorded session was
#      re-saved with the contributor list re-ordered (the email now listed
#      before "System" rather than after).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1 & notes: the 6 newly-lapsed sessions, one per sub-group, all dated
# 21/01/2026. Copy formatting from an existing "Not Recorded" row (row 3,
# style index 4) onto each header row, then update the status text.
# ---------------------------------------------------------------------------
$lapsedRows = @(183, 210, 237, 264, 291, 318)

$ws.Range("A3:I3").Copy()
foreach ($r in $lapsedRows) {
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $ws.Range("I" + $r).Value = "Not Recorded"
}

# ---------------------------------------------------------------------------
# 2. Class Statistics roll-up (K/L columns near the top of the sheet).
# ---------------------------------------------------------------------------
$ws.Range("L7").Value = 63   # Missing Sessions: 57 -> 63
$ws.Range("L8").Value = 6    # Pending Sessions: 12 -> 6

# ---------------------------------------------------------------------------
# 3. Group Statistics table: Missing (+1) / Pending (-1) for the six groups
#    whose 21/01/2026 session just lapsed.
# ---------------------------------------------------------------------------
$groupRows = @(21, 22, 23, 24, 25, 26)
foreach ($r in $groupRows) {
    $ws.Range("P" + $r).Value = $ws.Range("P" + $r).Value2 + 1
    $ws.Range("Q" + $r).Value = $ws.Range("Q" + $r).Value2 - 1
}

# ---------------------------------------------------------------------------
# 4. "Recorded By" re-ordering: "System, dnasr281@gmail.com"
#    -> "dnasr281@gmail.com, System" for every already-recorded session row.
# ---------------------------------------------------------------------------
$recordedByRows = @(8,9,10,12,14,15,17,18,23,34,35,36,38,40,41,43,44,49,60,61,62,64,66,67,69,70,75,86,87,88,90,92,93,95,96,101,112,113,114,116,118,119,121,122,127,138,139,140,142,144,145,147,148,153,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)

foreach ($r in $recordedByRows) {
    $ws.Range("G" + $r).Value = "dnasr281@gmail.com, System"
}
